$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# New task progress marker in C16 ("en proceso")
$ws.Range("C16").Value = "en proceso"

# New task added at the bottom of the list (row 30)
$ws.Range("A30").Value = "Validacion en creacion de cuota, no muestra los mensajes de error"

# Scroll the view down a bit and leave the new row selected, reflecting
# where the user ended up after adding the task.
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("B30").Select()
